# Update "想去人数" (want-to-go count) values on two sheets, reflecting
# refreshed data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1069
$ws1.Range("F8").Value = 209
$ws1.Range("F9").Value = 384
$ws1.Range("F10").Value = 646
$ws1.Range("F11").Value = 8
$ws1.Range("F12").Value = 501
$ws1.Range("F15").Value = 12498
$ws1.Range("F16").Value = 140

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1069
$ws4.Range("F10").Value = 209
$ws4.Range("F11").Value = 384
$ws4.Range("F12").Value = 647
$ws4.Range("F13").Value = 8
$ws4.Range("F14").Value = 501
$ws4.Range("F17").Value = 12498
$ws4.Range("F20").Value = 140

$wb.Save()
